# Auto-generated edit script: updates cryptos price/volume data
# per commit "Updated cryptos list on Mon Jul 22 23:52:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.576.74"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "'3.442.67"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'591.43"
$ws.Range("D6").Value = "'178.94"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "'3.439.86"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("D13").Value = "'4.042.95"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "'31.87"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "'0.132"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "'67.548.20"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D18").Value = "'3.441.68"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("D20").Value = "'13.94"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").Value = "'385.55"
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'71.28"
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "'0.530"
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("E33").Value = "  -7.62%  "
$ws.Range("D34").Value = "'23.45"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -7.41%  "
$ws.Range("D38").Value = "'161.22"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").Value = "'0.879"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("D42").Value = "'6.61"
$ws.Range("E42").Value = "  -7.51%  "
$ws.Range("D43").Value = "'4.52"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").Value = "'25.70"
$ws.Range("E44").Value = "  -5.59%  "
$ws.Range("D45").Value = "'0.0709"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").Value = "'25.83"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "'2.694.12"
$ws.Range("E47").Value = "  -6.51%  "
$ws.Range("D48").Value = "'41.20"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("D50").Value = "'324.84"
$ws.Range("E50").Value = "  -7.19%  "
$ws.Range("E51").Value = "  -4.96%  "
